# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 647
$wsOff.Range("C3").Value = 463
$wsOff.Range("D3").Value = 184
$wsOff.Range("E3").Value = 88
$wsOff.Range("G3").Value = 7

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 663
$wsDef.Range("C3").Value = 466
$wsDef.Range("D3").Value = 120
$wsDef.Range("E3").Value = 61
